$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update/add the review rows (names + dates) ---
$ws.Range("A2").Value = "Ellen"
$ws.Range("A3").Value = "Wilbert"
$ws.Range("A4").Value = "Josie"
$ws.Range("A5").Value = "Marcy"
$ws.Range("A6").Value = "Carmine"

# Copy the existing date formatting from B2 down to B3:B6 before
# writing the new values so every date cell shares the same style.
$ws.Range("B2").Copy()
$ws.Range("B3:B6").PasteSpecial(-4122)

$ws.Range("B2").Value = 43779
$ws.Range("B3").Value = 43779
$ws.Range("B4").Value = 43779
$ws.Range("B5").Value = 43779
$ws.Range("B6").Value = 43779

# --- Remove the old hyperlink on A2 ---
$ws.Range("A2").Hyperlinks.Delete()

# --- Apply the new "card" styling to the name column (A2:A6) ---
$rng = $ws.Range("A2:A6")
$rng.Font.Color = 0
$rng.Font.Name = "Docs-Calibri"
$rng.Interior.Color = 16777215
$rng.Borders.Weight = -4138
$rng.Borders.Color = 13421772
$rng.WrapText = $true

# --- Row heights for the header + data rows ---
$ws.Range("A1:B6").RowHeight = 15.75

# --- Selection moves to B2 ---
$null = $ws.Range("B2").Select()
